# Updates the cryptos price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.652.06"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "1.827.70"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'1.006"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'308.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").Value = "'0.4660"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.54%  "

$ws.Range("D8").Value = "'0.3605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "'0.07135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").Value = "'0.9042"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.32%  "

$ws.Range("D11").Value = "'0.07744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "1.820.01"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("D14").Value = "'5.260"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "'6.338"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").Value = "'87.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.23%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "'0.000008550"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").Value = "26.692.45"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("D21").Value = "'14.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("D22").Value = "'5.012"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'1.912"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("D25").Value = "'152.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "'17.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("D27").Value = "'1.976"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").Value = "'113.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "

$ws.Range("D29").Value = "'4.846"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").Value = "'0.08803"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").Value = "'3.148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.17%  "

$ws.Range("D32").Value = "'2.843"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("D33").Value = "'0.7383"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.02%  "

$ws.Range("D34").Value = "'1.158"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.88%  "

$ws.Range("D35").Value = "'4.445"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").Value = "'0.01924"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").Value = "'2.925"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").Value = "'0.05132"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").Value = "'6.877"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").Value = "'0.5070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").Value = "'0.1499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.69%  "

$ws.Range("E43").Value = "  +0.46%  "

$ws.Range("D44").Value = "'0.4673"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("D45").Value = "'1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").Value = "'9.988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").Value = "'98.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").Value = "'0.06063"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("D50").Value = "'63.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").Value = "'35.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
